$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.036.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.041.43'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.665'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.13'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.45%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.384'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0786'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('E11').Value = '  +1.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.80'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.338.62'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.829'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.74'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.035.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +29.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.999.51'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0896'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '236.98'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.44%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.49%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.125'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.77'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0627'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.54'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0893'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.70%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.22'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.24%  '
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.107'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.34'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.18'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +14.40%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.22'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +19.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0222'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.45'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.03'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.289.33'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.91'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.80'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.219.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.95'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.15%  '
